$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update individual values (totalEstadosMensuales / totalDeudores / totalSeguimientos adjustments)
$ws.Range("E5").Value = 1002

$ws.Range("C13").Value = 77
$ws.Range("E13").Value = 914

$ws.Range("D16").Value = 12
$ws.Range("E16").Value = 97

$ws.Range("E25").Value = 2919

$ws.Range("E41").Value = 239

$ws.Range("E46").Value = 498

$ws.Range("E48").Value = 92

$ws.Range("E55").Value = 404

$ws.Range("E58").Value = 273

$ws.Range("E67").Value = 565

$ws.Range("E90").Value = 479

$ws.Range("E99").Value = 294

$ws.Range("E103").Value = 2433

$ws.Range("E104").Value = 426

$ws.Range("E118").Value = 1175

$ws.Range("E119").Value = 64

$ws.Range("D122").Value = 50
$ws.Range("E122").Value = 531

$ws.Range("D123").Value = 80
$ws.Range("E123").Value = 1456

$ws.Range("E129").Value = 420

$ws.Range("E130").Value = 557

$ws.Range("E151").Value = 775

$ws.Range("E170").Value = 429

$ws.Range("E171").Value = 140

# Append a new row 180 with a client that migrated but is not in NOMOS
$ws.Range("A180").Value = "CONJUNTO MAZUREN SECTOR II AGRUPACION 3 - P.H."
$ws.Range("B180").Value = "830135053"
$ws.Range("C180").Value = 0
$ws.Range("D180").Value = 0
$ws.Range("E180").Value = 0
$ws.Range("F180").Value = "NO ESTA EN NOMOS"
$ws.Range("G180").Value = ""
